# Auto-generated edit script applying the diff's numeric value changes
# to Sheet1 of the Betfair odds workbook (row r=2..11, various columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "F2" = 2.88
    "G2" = 3
    "J2" = 3.75
    "L2" = 1.3
    "M2" = 1.02
    "P2" = 2.34
    "Q2" = 1.6
    "R2" = 1.54
    "S2" = 2.36
    "U2" = 2.48
    "X2" = 28
    "AB2" = 17
    "AE2" = 25
    "AF2" = 24
    "AK2" = 1000
    "AO2" = 17
    "H3" = 1.39
    "I3" = 1.5
    "J3" = 4.1
    "Q3" = 1.72
    "G4" = 1.86
    "H4" = 4.7
    "J4" = 3.35
    "Q4" = 2.08
    "H5" = 1.48
    "J5" = 3.6
    "I6" = 1.71
    "J6" = 3.7
    "F8" = 1.49
    "G8" = 1.5
    "H8" = 7.2
    "I8" = 7.4
    "J8" = 5.2
    "K8" = 5.3
    "Q8" = 1.59
    "R8" = 1.65
    "AE8" = 95
    "AI8" = 960
    "AL8" = 29
    "AM8" = 90
    "AN8" = 5.5
    "AO8" = 100
    "F10" = 2.04
    "H10" = 3.7
    "K10" = 3.8
    "P10" = 1.62
    "Q10" = 1.77
    "R10" = 1.29
    "U10" = 1.99
    "W10" = 1.83
    "X10" = 15.5
    "Y10" = 18.5
    "AC10" = 10.5
    "AD10" = 23
    "AF10" = 18.5
    "AG10" = 15
    "AH10" = 1000
    "AM10" = 1000
    "AN10" = 1000
    "F11" = 1.9
    "G11" = 2.12
    "I11" = 5.4
    "J11" = 3.05
    "K11" = 4.5
    "P11" = 1.82
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

